$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339, shifting existing rows 339:427 down to 340:428
$ws.Rows.Item(339).Insert()

# Populate the newly inserted row 339 with the new data
$ws.Cells.Item(339, 1).Value = 4
$ws.Cells.Item(339, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(339, 3).Value = "Los Lagos"
$ws.Cells.Item(339, 4).Value = 44964
$ws.Cells.Item(339, 5).Value = 10
$ws.Cells.Item(339, 6).Value = 100112045
$ws.Cells.Item(339, 7).Value = "Zapallo"
$ws.Cells.Item(339, 8).Value = "Paine"
$ws.Cells.Item(339, 9).Value = "1a (cosecha)"
$ws.Cells.Item(339, 10).Value = 1200
$ws.Cells.Item(339, 11).Value = 550
$ws.Cells.Item(339, 12).Value = 600
$ws.Cells.Item(339, 13).Value = 575
$ws.Cells.Item(339, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(339, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(339, 16).Value = 575
$ws.Cells.Item(339, 17).Value = 1
$ws.Cells.Item(339, 18).Value = "Hortaliza"
